$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 32, shifting existing rows 32:35 down to 33:36.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new data record.
$ws.Cells.Item(32, 1).Value = 10
$ws.Cells.Item(32, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(32, 3).Value = "La Araucanía"
$ws.Cells.Item(32, 4).Value = 44783
$ws.Cells.Item(32, 5).Value = 9
$ws.Cells.Item(32, 6).Value = 100112042
$ws.Cells.Item(32, 7).Value = "Locoto"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 90
$ws.Cells.Item(32, 11).Value = 2700
$ws.Cells.Item(32, 12).Value = 2700
$ws.Cells.Item(32, 13).Value = 2700
$ws.Cells.Item(32, 14).Value = "$/kilo"
$ws.Cells.Item(32, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(32, 16).Value = 2700
$ws.Cells.Item(32, 17).Value = 1
$ws.Cells.Item(32, 18).Value = "Hortaliza"
